# TC08_Canine_StudyUBC02-Breed_Diagnosis_PrimDiseaseSite.xlsx
#
# The "CasesTab" Cypher query (cell B2 on the "startup" sheet) dropped its
# trailing `co.cohort_description AS Cohort` projection - the cohort column
# is no longer being returned, so the RETURN clause now ends at
# `Response to Treatment`.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newQuery = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN [''UBC02''] and demo.breed in [''Mixed Breed'', ''Scottish Terrier'',''Shetland Sheepdog'']and diag.disease_term in [''Bladder Cancer'',''Healthy Control''] and diag.primary_disease_site in [''Bladder'', ''Bladder, Urethra, Vagina'']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '''') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '''') AS `Study Code` ,
        coalesce(s.clinical_study_type, '''') AS  `Study Type`,
        coalesce(demo.breed, '''') AS Breed ,
        coalesce(diag.disease_term, '''') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '''') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '''') AS Age ,
        coalesce(demo.sex, '''') AS Sex ,
        coalesce(demo.neutered_indicator, '''') AS `Neutered Status`,
        coalesce(demo.weight, '''') AS `Weight (kg)`,
        coalesce(diag.best_response, '''') AS `Response to Treatment`'

$ws.Range("B2").Value = $newQuery

# Match the author's saved cursor position on that same cell.
$ws.Activate()
$ws.Range("B2").Select()
